# Updates Poroto granado (Feria Lagunitas de Puerto Montt) weekly price data.
# Refreshes Fecha (D), Volumen (J), Precio minimo/maximo/promedio ponderado (K/L/M),
# Origen (O) and Precio $/Kg (P) for rows 3-21 (row 16 unchanged) per the weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44243
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 33000
$ws.Range("L3").Value = 33000
$ws.Range("M3").Value = 33000
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 1320

$ws.Range("D4").Value = 44203
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 30000
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 1200

$ws.Range("D5").Value = 44568
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 35000
$ws.Range("L5").Value = 35000
$ws.Range("M5").Value = 35000
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 1400

$ws.Range("D6").Value = 44222
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 40000
$ws.Range("L6").Value = 40000
$ws.Range("M6").Value = 40000
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 1600

$ws.Range("D7").Value = 44232
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 40000
$ws.Range("L7").Value = 40000
$ws.Range("M7").Value = 40000
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 1600

$ws.Range("D8").Value = 44253
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 30000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 30000
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 1200

$ws.Range("D9").Value = 44246
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 31000
$ws.Range("L9").Value = 31000
$ws.Range("M9").Value = 31000
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1240

$ws.Range("D10").Value = 44572
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 35000
$ws.Range("L10").Value = 35000
$ws.Range("M10").Value = 35000
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 1400

$ws.Range("D11").Value = 44202
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 30000
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1200

$ws.Range("D12").Value = 44271
$ws.Range("J12").Value = 40
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 30000
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1200

$ws.Range("D13").Value = 44201
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 30000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 30000
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1200

$ws.Range("D14").Value = 44215
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 35000
$ws.Range("L14").Value = 35000
$ws.Range("M14").Value = 35000
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 1400

$ws.Range("D15").Value = 44204
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 37000
$ws.Range("L15").Value = 37000
$ws.Range("M15").Value = 37000
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1480

$ws.Range("D17").Value = 44211
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 42000
$ws.Range("L17").Value = 42000
$ws.Range("M17").Value = 42000
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 1680

$ws.Range("D18").Value = 44239
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 35000
$ws.Range("L18").Value = 35000
$ws.Range("M18").Value = 35000
$ws.Range("O18").Value = "Región del Maule"
$ws.Range("P18").Value = 1400

$ws.Range("D19").Value = 44250
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = 30000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 30000
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1200

$ws.Range("D20").Value = 44236
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 38000
$ws.Range("L20").Value = 38000
$ws.Range("M20").Value = 38000
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 1520

$ws.Range("D21").Value = 44264
$ws.Range("J21").Value = 30
$ws.Range("K21").Value = 29000
$ws.Range("L21").Value = 29000
$ws.Range("M21").Value = 29000
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1160

